$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 16159
$ws1.Range("F24").Value = 6474
$ws1.Range("F26").Value = 9
$ws1.Range("F28").Value = 4
$ws1.Range("F32").Value = 165
$ws1.Range("F33").Value = 4705
$ws1.Range("F34").Value = 14

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 16159
$ws4.Range("F25").Value = 6474
$ws4.Range("F27").Value = 9
$ws4.Range("F29").Value = 4
$ws4.Range("F34").Value = 165
$ws4.Range("F35").Value = 4705
$ws4.Range("F36").Value = 14
